# BusinessSupportDatasheet.pptx — translation-pipeline re-save.
#  1) The four "datetimeFigureOut" placeholders (notes master, slide
#     master, and both slide layouts) get re-cached from 10/11/21 to
#     11/18/21, as PowerPoint does whenever such a deck is reopened and
#     resaved on a later date.
#  2) The slide-1 title run "Adobe 支持产品/服务" is split in two so the
#     translated tail ("支持计划") can carry its own language tag.

$p = $ppt.ActivePresentation

# --- 1a. Notes master date placeholder ----------------------------------
# (the notes master's Shapes collection is not writable in this host, so
# go through the HeadersFooters façade instead, same as the Header/Footer
# dialog would)
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "11/18/21"

# --- 1b. Slide master date placeholder (shape 4) -----------------------
$master = $p.SlideMaster
$master.Shapes.Item(4).TextFrame.TextRange.Text = "11/18/21"

# --- 1c/1d. Both slide layout date placeholders -------------------------
$layout1 = $master.CustomLayouts.Item(1)
$layout1.Shapes.Item(3).TextFrame.TextRange.Text = "11/18/21"

$layout2 = $master.CustomLayouts.Item(2)
$layout2.Shapes.Item(1).TextFrame.TextRange.Text = "11/18/21"

# --- 2. Slide 1 title: split the run and retag the translated half -----
$slide1 = $p.Slides.Item(1)
$title = $slide1.Shapes.Item(3)
$tr = $title.TextFrame.TextRange

# Build the text with the Chinese tail first, while it is the only run,
# so that the language tag lands on that run; then prepend "Adobe " and
# restore that leading run's own (untouched) language.
$tr.Text = "支持计划"
$tr.LanguageID = "ja-JP"
$tr.InsertBefore("Adobe ")

$lead = $tr.Characters(1, 6)
$lead.LanguageID = "en-US"
